# "Listener and Log4j, allure report added"
#
# The underlying data change is: the plaintext password stored in
# Sheet1!B2/B3 (and cached in the shared string table) is masked out,
# while the original password text is preserved as the cached display
# text of the two "Password" hyperlinks, and the current cell selection
# is moved from C4 to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$plainPassword = "Suraj@123"
$maskedPassword = "************"

# Re-create the existing mailto hyperlinks on B2/B3 so their cached
# "display" text keeps showing the real password, exactly like Excel
# does when the visible cell text and the hyperlink's own display text
# diverge.
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:" + $plainPassword, [System.Type]::Missing, [System.Type]::Missing, $plainPassword) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:" + $plainPassword, [System.Type]::Missing, [System.Type]::Missing, $plainPassword) | Out-Null

# Mask the password that is actually shown in the cells.
$ws.Range("B2").Value = $maskedPassword
$ws.Range("B3").Value = $maskedPassword

# Move the active selection from C4 to B2.
$ws.Range("B2").Select()
